$d = $word.ActiveDocument

# 1. Scatter plot bullet: add "vs bedrooms" comparison text.
$d.Content.Find.Execute(
    "Created a scatter plot to analyze the relationship between sold_price and zipcode for identifying high-demand areas and outliers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created a scatter plot to analyze the relationship between sold_price vs zipcode vs bedrooms for identifying high-demand areas and outliers.",
    2)

# 2. Box plot bullet -> bar chart bullet about number of bedrooms.
$d.Content.Find.Execute(
    "Generated a box plot to examine property trends based on construction year (year_built) and sales (sold_price).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Generated a bar chart to examine number of properties and number of bedrooms properties to know about which property has how many bedrooms.",
    2)

# 3. Add a new bullet after "Saved the cleaned dataset..." describing the final dataset shape.
$saved = $d.Paragraphs(20)
$saved.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(21)
$newPara.Range.Text = "Final clean dataset has 4370 observations and 16 features."
